# Helper: replace a short, previously-inserted placeholder token with a
# literal piece of text (used for characters such as straight quotes /
# apostrophes that Word's Find.Execute "smart quotes" autocorrect would
# otherwise mangle when passed as a Find.Execute replacement string).
function Replace-Token($doc, $paraIndex, $token, $replacement) {
    $searchRange = $doc.Paragraphs($paraIndex).Range.Duplicate
    $found = $searchRange.Find.Execute($token)
    if ($found) {
        $searchRange.Text = $replacement
    }
    return $found
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert the new "Meta description: ..." paragraph right after the
#    title (Heading1) paragraph.
# ---------------------------------------------------------------------

$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

# Borrow the run layout (leading empty run + bold run) from the
# duplicate bold paragraph further down the document so the XML shape
# matches exactly (<w:r/><w:r><w:rPr><w:b/></w:rPr>...).
$paraCountBefore = $d.Paragraphs.Count
$dupBoldPara = $d.Paragraphs($paraCountBefore - 1)
$metaPara.Range.FormattedText = $dupBoldPara.Range.FormattedText

# Turn the copied bold text into "Meta description".
$metaRange = $d.Paragraphs(2).Range
$metaRange.Find.Execute(
    "Play 9 Tigers Slot for Free - Review of Features & Symbols",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Meta description", 2) | Out-Null

# Append the (non-bold) remainder of the sentence right after it.
$metaRange2 = $d.Paragraphs(2).Range
$insertPoint = $d.Range($metaRange2.End - 1, $metaRange2.End - 1)
$insertPoint.InsertAfter(": Read our 9 Tigers slot review and play for free! Find out about the game's features, symbols, and potential payouts.")

# ---------------------------------------------------------------------
# 2) Remove the now-duplicated bold "Play 9 Tigers Slot..." paragraph
#    near the end of the document, and rewrite the italic paragraph's
#    text with the new AI image-prompt copy.
# ---------------------------------------------------------------------

$paraCount = $d.Paragraphs.Count
$dupPara = $d.Paragraphs($paraCount - 1)
$dupPara.Range.Delete()

$paraCount = $d.Paragraphs.Count
$italicParaIndex = $paraCount
$italicPara = $d.Paragraphs($italicParaIndex)
$italicRange = $italicPara.Range
$textOnly = $d.Range($italicRange.Start, $italicRange.End - 1)

# Use placeholder tokens in place of straight quotes / apostrophe so the
# bulk Find.Execute replacement is not mangled by smart-quote
# autocorrection, then patch the placeholders back in afterwards via
# small in-place Range.Text assignments (which do not autocorrect).
$placeholder = 'Create a cartoon-style feature image for ZQUOTEZ9 TigersZQUOTEZ game that showcases a happy Maya warrior with glasses. The image should be colorful and vibrant, with the warrior holding up a winning combination of Fire and Water Tigers on the reels. In the background, elements of Yin and Yang can be shown. The warrior can be dressed in traditional Maya clothing, with a big smile on their face, and glasses which make them look intelligent and fun-loving. Overall, the image should reflect the gameZAPOSZs Asian-inspired theme and the joy of winning.'

$textOnly.Find.Execute(
    "Read our 9 Tigers slot review and play for free! Find out about the game's features, symbols, and potential payouts.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $placeholder, 2) | Out-Null

Replace-Token $d $italicParaIndex "ZQUOTEZ" '"' | Out-Null
Replace-Token $d $italicParaIndex "ZQUOTEZ" '"' | Out-Null
Replace-Token $d $italicParaIndex "ZAPOSZ" "'" | Out-Null

Write-Output "done"
